$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 327
$ws1.Range("F3").Value = 90
$ws1.Range("F4").Value = 488
$ws1.Range("F5").Value = 4827
$ws1.Range("F9").Value = 738
$ws1.Range("F10").Value = 219

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 327
$ws4.Range("F3").Value = 90
$ws4.Range("F4").Value = 488
$ws4.Range("F5").Value = 4827
$ws4.Range("F9").Value = 738
$ws4.Range("F11").Value = 219
